$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C10) from serial 45233 to 45243
$ws.Range("C2:C10").Value = 45243
